$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1765.1111
$ws.Range("I70").Value = 1580.4736
$ws.Range("J70").Value = 1971.4706
$ws.Range("K70").Value = 4741.4208
$ws.Range("L70").Value = 5914.4118
$ws.Range("M70").Value = -4471.4208
$ws.Range("N70").Value = -6454.4118
$ws.Range("H73").Value = 1765.1111
$ws.Range("I73").Value = 1580.4736
$ws.Range("J73").Value = 1971.4706
$ws.Range("K73").Value = 4741.4208
$ws.Range("L73").Value = 5914.4118
$ws.Range("M73").Value = -3805.4208
$ws.Range("N73").Value = -7786.4118
$ws.Range("H108").Value = 48033.332
$ws.Range("J108").Value = 48033.332
$ws.Range("L108").Value = 48033.332
$ws.Range("N108").Value = -55713.332
$ws.Range("H113").Value = 2584.5
$ws.Range("I113").Value = 2343.1667
$ws.Range("J113").Value = 2765.5
$ws.Range("K113").Value = 2343.1667
$ws.Range("L113").Value = 2765.5
$ws.Range("M113").Value = 910.8332999999998
$ws.Range("N113").Value = -9273.5
$ws.Range("H121").Value = 1173.3077
$ws.Range("J121").Value = 1037.75
$ws.Range("L121").Value = 3113.25
$ws.Range("N121").Value = -6607.25
$ws.Range("H137").Value = 903181.5600000001
$ws.Range("I137").Value = 2979.3
$ws.Range("J137").Value = 1962243
$ws.Range("K137").Value = 8937.900000000001
$ws.Range("L137").Value = 5886729
$ws.Range("M137").Value = -6387.900000000001
$ws.Range("N137").Value = -5891829

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1073.64
$ws.Range("I2").Value = 1097.4348
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 1097.4348
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = -984.4348
$ws.Range("N2").Value = -1026
$ws.Range("H32").Value = 26823.559
$ws.Range("I32").Value = 36639.566
$ws.Range("J32").Value = 4171.231
$ws.Range("K32").Value = 36639.566
$ws.Range("L32").Value = 4171.231
$ws.Range("M32").Value = -36352.566
$ws.Range("N32").Value = -4745.231
$ws.Range("H45").Value = 1640.5758
$ws.Range("I45").Value = 1750.4584
$ws.Range("J45").Value = 1347.5555
$ws.Range("K45").Value = 1750.4584
$ws.Range("L45").Value = 1347.5555
$ws.Range("M45").Value = -1373.4584
$ws.Range("N45").Value = -2101.5555
$ws.Range("H63").Value = 2300.5557
$ws.Range("I63").Value = 2300.5557
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2300.5557
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1614.5557
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 2300.5557
$ws.Range("I66").Value = 2300.5557
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 11502.7785
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -8070.7785
$ws.Range("N66").ClearContents()
$ws.Range("H116").Value = 1073.64
$ws.Range("I116").Value = 1097.4348
$ws.Range("J116").Value = 800
$ws.Range("K116").Value = 1097.4348
$ws.Range("L116").Value = 800
$ws.Range("M116").Value = 1196.5652
$ws.Range("N116").Value = -5388

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1073.64
$ws.Range("I3").Value = 1097.4348
$ws.Range("J3").Value = 800
$ws.Range("K3").Value = 1097.4348
$ws.Range("L3").Value = 800
$ws.Range("M3").Value = -983.4348
$ws.Range("N3").Value = -1028
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H120").Value = 50000
$ws.Range("J120").Value = 50000
$ws.Range("L120").Value = 50000
$ws.Range("N120").Value = -59676

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 702297.1
$ws.Range("I31").Value = 5521.727
$ws.Range("J31").Value = 1196782.9
$ws.Range("K31").Value = 5521.727
$ws.Range("L31").Value = 1196782.9
$ws.Range("M31").Value = -5226.727
$ws.Range("N31").Value = -1197372.9
$ws.Range("H34").Value = 702297.1
$ws.Range("I34").Value = 5521.727
$ws.Range("J34").Value = 1196782.9
$ws.Range("K34").Value = 5521.727
$ws.Range("L34").Value = 1196782.9
$ws.Range("M34").Value = -5319.727
$ws.Range("N34").Value = -1197186.9
$ws.Range("H41").Value = 59999
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H51").Value = 22847.5
$ws.Range("I51").Value = 7090
$ws.Range("J51").Value = 28100
$ws.Range("K51").Value = 7090
$ws.Range("L51").Value = 28100
$ws.Range("M51").Value = -6354
$ws.Range("N51").Value = -29572
$ws.Range("H60").Value = 16098.25
$ws.Range("J60").Value = 18100
$ws.Range("L60").Value = 18100
$ws.Range("N60").Value = -19122
$ws.Range("H61").Value = 22847.5
$ws.Range("I61").Value = 7090
$ws.Range("J61").Value = 28100
$ws.Range("K61").Value = 7090
$ws.Range("L61").Value = 28100
$ws.Range("M61").Value = -6742
$ws.Range("N61").Value = -28796
$ws.Range("H99").Value = 1294.625
$ws.Range("J99").Value = 1228.5
$ws.Range("L99").Value = 1228.5
$ws.Range("N99").Value = -4224.5
$ws.Range("H126").Value = 1294.625
$ws.Range("J126").Value = 1228.5
$ws.Range("L126").Value = 3685.5
$ws.Range("N126").Value = -8625.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2042.5625
$ws.Range("J34").Value = 2448.6155
$ws.Range("L34").Value = 7345.8465
$ws.Range("N34").Value = -7513.8465
$ws.Range("H51").Value = 1892.6666
$ws.Range("I51").Value = 1839.75
$ws.Range("J51").Value = 1998.5
$ws.Range("K51").Value = 5519.25
$ws.Range("L51").Value = 5995.5
$ws.Range("M51").Value = -5059.25
$ws.Range("N51").Value = -6915.5
$ws.Range("H55").Value = 7756.533
$ws.Range("I55").Value = 1933.3334
$ws.Range("J55").Value = 9212.333000000001
$ws.Range("K55").Value = 5800.0002
$ws.Range("L55").Value = 27636.999
$ws.Range("M55").Value = -5623.0002
$ws.Range("N55").Value = -27990.999
$ws.Range("H64").Value = 3666.6667
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 15000
$ws.Range("M64").Value = -2730
$ws.Range("N64").Value = -15540
$ws.Range("H67").Value = 3666.6667
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 15000
$ws.Range("M67").Value = -2064
$ws.Range("N67").Value = -16872
$ws.Range("H107").Value = 622.61224
$ws.Range("I107").Value = 286.23254
$ws.Range("J107").Value = 3033.3333
$ws.Range("K107").Value = 858.6976199999999
$ws.Range("L107").Value = 9099.999899999999
$ws.Range("M107").Value = 1061.30238
$ws.Range("N107").Value = -12939.9999
$ws.Range("H121").Value = 750
$ws.Range("I121").Value = 1000
$ws.Range("J121").Value = 666.6667
$ws.Range("K121").Value = 3000
$ws.Range("L121").Value = 2000.0001
$ws.Range("M121").Value = -1690
$ws.Range("N121").Value = -4620.0001
$ws.Range("H131").Value = 1016.11664
$ws.Range("I131").Value = 339.29413
$ws.Range("J131").Value = 1283.6976
$ws.Range("K131").Value = 1017.88239
$ws.Range("L131").Value = 3851.0928
$ws.Range("M131").Value = 4022.11761
$ws.Range("N131").Value = -13931.0928

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 20000
$ws.Range("J58").Value = 20000
$ws.Range("L58").Value = 20000
$ws.Range("N58").Value = -20554

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2556.913
$ws.Range("I7").Value = 2288.1333
$ws.Range("J7").Value = 3060.875
$ws.Range("K7").Value = 2288.1333
$ws.Range("L7").Value = 3060.875
$ws.Range("M7").Value = -2176.1333
$ws.Range("N7").Value = -3284.875
$ws.Range("H40").Value = 2930.375
$ws.Range("I40").Value = 2795.5908
$ws.Range("J40").Value = 3226.9
$ws.Range("K40").Value = 2795.5908
$ws.Range("L40").Value = 3226.9
$ws.Range("M40").Value = -2659.5908
$ws.Range("N40").Value = -3498.9
$ws.Range("H45").Value = 16825.715
$ws.Range("I45").Value = 9000
$ws.Range("J45").Value = 18130
$ws.Range("K45").Value = 9000
$ws.Range("L45").Value = 18130
$ws.Range("M45").Value = -8593
$ws.Range("N45").Value = -18944
$ws.Range("H46").Value = 836.4286
$ws.Range("I46").Value = 647.36365
$ws.Range("K46").Value = 647.36365
$ws.Range("M46").Value = -459.36365
$ws.Range("H68").Value = 3450
$ws.Range("I68").Value = 2600
$ws.Range("J68").Value = 3960
$ws.Range("K68").Value = 2600
$ws.Range("L68").Value = 3960
$ws.Range("M68").Value = -1851
$ws.Range("N68").Value = -5458
$ws.Range("H71").Value = 3450
$ws.Range("I71").Value = 2600
$ws.Range("J71").Value = 3960
$ws.Range("K71").Value = 13000
$ws.Range("L71").Value = 19800
$ws.Range("M71").Value = -9256
$ws.Range("N71").Value = -27288
$ws.Range("H126").Value = 2556.913
$ws.Range("I126").Value = 2288.1333
$ws.Range("J126").Value = 3060.875
$ws.Range("K126").Value = 6864.3999
$ws.Range("L126").Value = 9182.625
$ws.Range("M126").Value = -4394.3999
$ws.Range("N126").Value = -14122.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3668.6667
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376
$ws.Range("H65").Value = 3668.6667
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880
$ws.Range("H86").Value = 20671.428
$ws.Range("J86").Value = 20671.428
$ws.Range("L86").Value = 20671.428
$ws.Range("N86").Value = -22917.428
$ws.Range("H89").Value = 20671.428
$ws.Range("J89").Value = 20671.428
$ws.Range("L89").Value = 103357.14
$ws.Range("N89").Value = -114589.14
$ws.Range("H107").Value = 1109.8462
$ws.Range("I107").Value = 678.64703
$ws.Range("J107").Value = 1924.3334
$ws.Range("K107").Value = 2035.94109
$ws.Range("L107").Value = 5773.0002
$ws.Range("M107").Value = -115.9410899999998
$ws.Range("N107").Value = -9613.0002
